$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update manualStatus column (I) for rows 25, 26, 27, 29 from plain
# numeric values to bracketed-list text values.
$ws.Range("I25").Value = "[8]"
$ws.Range("I26").Value = "[4]"
$ws.Range("I27").Value = "[4,256]"
$ws.Range("I29").Value = "[4]"

# Widen the fastqFileName column (F) so the long file names are readable.
$ws.Columns.Item(6).ColumnWidth = 51.6

# Move the active selection to I30.
$ws.Range("I30").Select()
